$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q6)
$ws.Range("B2").Value = 0.3446948460339797
$ws.Range("C2").Value = 0.3446948460339797
$ws.Range("D2").Value = 0.1445232426843725
$ws.Range("E2").Value = 0.3801621268411315
$ws.Range("F2").Value = 0.1731862872813958
$ws.Range("G2").Value = 7

# Row 3 (Q7)
$ws.Range("B3").Value = 0.3269786427793283
$ws.Range("C3").Value = 0.3312344728690593
$ws.Range("D3").Value = 0.1381123416715895
$ws.Range("E3").Value = 0.3716346884664959
$ws.Range("F3").Value = 0.1852485889866795
$ws.Range("G3").Value = 11

# Row 4 (Q8)
$ws.Range("B4").Value = 0.2602905224375459
$ws.Range("C4").Value = 0.2850643385747676
$ws.Range("D4").Value = 0.1040594201288039
$ws.Range("E4").Value = 0.3225824237753879
$ws.Range("F4").Value = 0.200854463779107
$ws.Range("G4").Value = 10

# Row 5 (Q9)
$ws.Range("B5").Value = 0.3235260940633222
$ws.Range("C5").Value = 0.3235260940633222
$ws.Range("D5").Value = 0.1350471282558087
$ws.Range("E5").Value = 0.3674875892541253
$ws.Range("F5").Value = 0.186326885464196
$ws.Range("G5").Value = 8

# Row 6
$ws.Range("B6").Value = 0.3105992790096791
$ws.Range("C6").Value = 0.3148849937619054
$ws.Range("D6").Value = 0.1282073827774651
$ws.Range("E6").Value = 0.358060585344806
$ws.Range("F6").Value = 0.190444504571002
$ws.Range("G6").Value = 8

# Row 7
$ws.Range("B7").Value = 0.2715293295939802
$ws.Range("C7").Value = 0.2808173377963946
$ws.Range("D7").Value = 0.1058114480826604
$ws.Range("E7").Value = 0.325286716732578
$ws.Range("F7").Value = 0.1934695922577707
$ws.Range("G7").Value = 7

# Row 8
$ws.Range("B8").Value = 0.29537703190916
$ws.Range("C8").Value = 0.29537703190916
$ws.Range("D8").Value = 0.119218943402999
$ws.Range("E8").Value = 0.3452809629895616
$ws.Range("F8").Value = 0.1958714448515681
$ws.Range("G8").Value = 6

# Row 9
$ws.Range("B9").Value = 0.3108436987415587
$ws.Range("C9").Value = 0.3108436987415587
$ws.Range("D9").Value = 0.1261775239899553
$ws.Range("E9").Value = 0.3552147575621758
$ws.Range("F9").Value = 0.1883201070813918
$ws.Range("G9").Value = 6

# Row 10
$ws.Range("B10").Value = 0.2851881662423609
$ws.Range("C10").Value = 0.2851881662423609
$ws.Range("D10").Value = 0.1084188778872009
$ws.Range("E10").Value = 0.3292702201645344
$ws.Range("F10").Value = 0.1840060723268408
$ws.Range("G10").Value = 5

# Row 11
$ws.Range("B11").Value = 0.2443419770806189
$ws.Range("C11").Value = 0.2443419770806189
$ws.Range("D11").Value = 0.0744009712500677
$ws.Range("E11").Value = 0.2727654143216616
$ws.Range("F11").Value = 0.1399903305537064
$ws.Range("G11").Value = 4

$wb.Save()
